$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.459.28'
$ws.Range("E2").Value = '  -0.26%  '

$ws.Range("D3").Value = '1.581.45'
$ws.Range("E3").Value = '  -0.62%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.47'
$ws.Range("E5").Value = '  +0.27%  '

$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("E7").Value = '  +0.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.29'
$ws.Range("E8").Value = '  +0.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.03'
$ws.Range("E9").Value = '  -1.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.247'
$ws.Range("E10").Value = '  -1.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0893'
$ws.Range("E12").Value = '  +0.75%  '

$ws.Range("D13").Value = '1.808.77'
$ws.Range("E13").Value = '  -0.56%  '

$ws.Range("D14").Value = '1.583.67'
$ws.Range("E14").Value = '  -0.46%  '

$ws.Range("E15").Value = '  -1.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.520'
$ws.Range("E16").Value = '  -1.84%  '

$ws.Range("D17").Value = '28.483.93'
$ws.Range("E17").Value = '  -0.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.20'
$ws.Range("E18").Value = '  -1.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '230.63'
$ws.Range("E19").Value = '  -0.92%  '

$ws.Range("E20").Value = '  -0.85%  '

$ws.Range("E21").Value = '  -2.20%  '

$ws.Range("E22").Value = '  +0.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.90'
$ws.Range("E23").Value = '  -3.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.15'
$ws.Range("E24").Value = '  -2.03%  '

$ws.Range("E25").Value = '  +4.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.74'
$ws.Range("E26").Value = '  +0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.03'
$ws.Range("E27").Value = '  -1.61%  '

$ws.Range("E28").Value = '  -1.80%  '

$ws.Range("E29").Value = '  -2.30%  '

$ws.Range("E30").Value = '  +0.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0483'
$ws.Range("E31").Value = '  +2.49%  '

$ws.Range("E32").Value = '  -1.59%  '

$ws.Range("E33").Value = '  -1.34%  '

$ws.Range("E34").Value = '  -2.19%  '

$ws.Range("D35").Value = '1.398.65'
$ws.Range("E35").Value = '  -0.36%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.09'
$ws.Range("E36").Value = '  +7.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.52'
$ws.Range("E37").Value = '  -4.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.36'
$ws.Range("E38").Value = '  +0.57%  '

$ws.Range("E39").Value = '  +0.81%  '

$ws.Range("E40").Value = '  -0.61%  '

$ws.Range("E41").Value = '  -3.91%  '

$ws.Range("E42").Value = '  +0.24%  '

$ws.Range("E43").Value = '  +1.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.790'
$ws.Range("E44").Value = '  -2.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0466'
$ws.Range("E45").Value = '  -0.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.46'
$ws.Range("E46").Value = '  -3.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.960'
$ws.Range("E47").Value = '  -2.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.40'
$ws.Range("E48").Value = '  +0.48%  '

$ws.Range("D49").Value = '1.720.17'
$ws.Range("E49").Value = '  -0.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.72'
$ws.Range("E50").Value = '  -0.74%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0102'
$ws.Range("E51").Value = '  -2.70%  '
